# Quarterly indexing esoteric bug-fix operation
#
# Column A (rows 2-73) holds per-row "as of" dates that were meant to land on
# the 15th of the month AFTER the originally stored 1st-of-month date, but
# were off by one reporting cycle. Shift every date in column A (rows 2-73)
# forward by one month, pinned to the 15th of that new month.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 73; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $oldVal = $cell.Value()
    $oldDate = [datetime]$oldVal

    $shifted = $oldDate.AddMonths(1)
    $newDate = Get-Date -Year $shifted.Year -Month $shifted.Month -Day 15
    $newDate = $newDate.Date

    $cell.Value = $newDate
}
